$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old hack entries in rows 3-10, keeping only the header (row 1)
# and the single remaining data row (row 2).
$ws.Range("A3:I10").EntireRow.Delete()

# Row 2 now describes a new/placeholder entry: most fields are cleared,
# "Important" stays checked, and the Url points at the new custom hack.
$ws.Range("A2:G2").ClearContents()
$ws.Range("H2").Value = $true
$ws.Range("I2").Value = "http://www.arcadecollecting.com/hacks/bagman/"
